# Update data parsing logic
# Appends one new data row (row 86) to each of the 4 worksheets,
# mirroring the newest sample captured by the logger.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newRow = 86

# Per-sheet new-row payloads: A (serial date), B, C, D, E (text), F, G, H, I (numbers)
$sheetsData = @(
    @{
        Index = 1
        A = 45872.45982638889
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x28"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 296
        I = 7
    },
    @{
        Index = 2
        A = 45872.45982638889
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x2C"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 300
        I = 25
    },
    @{
        Index = 3
        A = 45872.45982638889
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5F"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 95
        I = 15
    },
    @{
        Index = 4
        A = 45872.45982638889
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x76"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 118
        I = 9
    }
)

foreach ($entry in $sheetsData) {
    $ws = $wb.Worksheets.Item($entry.Index)

    # Column A keeps the same date/time number format used by the rows above it.
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($newRow, 1).Value = $entry.A

    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E

    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = [double]$entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}
